# Update the CDA Logical Model "url" StructureDefinition metadata sheet
# to the ST.r2b release (2.0.1-sd-202510-matchbox-patch / 2025-10-29) and
# add the new "Jurisdiction" metadata row (FHIR R5 StructureDefinition
# gained a jurisdiction element that the IG publisher now emits, with an
# empty value for this resource).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1. Bump the "Version" value (row 3, column B).
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# 2. Bump the "Date" value (row 8, column B).
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# 3. Insert a new "Jurisdiction" row right after "Contact" (row 10) and
#    before "Description" (old row 11), pushing everything else down by one.
$ws.Rows.Item(11).Insert()

# Carry the existing row-body formatting (border/alignment) down onto the
# freshly inserted row instead of leaving Excel's bare default style.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
